$wb = $excel.ActiveWorkbook

# Duplicate the "Italy" sheet to create the new "Spain" sheet, placed right after Italy.
$italy = $wb.Worksheets.Item("Italy")
$italy.Copy($null, $italy)
$newSheet = $wb.Worksheets.Item($italy.Index + 1)
$newSheet.Name = "Spain"

# Update the market-specific cell contents. Setting B4 before B2 keeps the shared-string
# table ordering (NGC-3103/T2044 then Spain Market) consistent with the target workbook.
$newSheet.Range("B4").Value = "NGC-3103/T2044"
$newSheet.Range("B2").Value = "Spain Market"

# The "notes" column (D) rows grow a bit taller on the new sheet.
$newSheet.Rows.Item(3).RowHeight = 28.8
$newSheet.Rows.Item(4).RowHeight = 28.8
$newSheet.Rows.Item(5).RowHeight = 28.8

# Columns B-D shrink to fit the shorter Spain content (engine persists column widths
# ~0.8333 characters wider than requested, so the inputs below are pre-compensated to
# land as close as possible to the real target widths of 15.21875 / 10.6640625 / 19.6640625).
$newSheet.Columns.Item(2).ColumnWidth = 14.333333333333334
$newSheet.Columns.Item(3).ColumnWidth = 9.833333333333334
$newSheet.Columns.Item(4).ColumnWidth = 18.833333333333332

# Restore the Italy sheet's selection (no longer the active tab) to the full data range.
$italy.Activate()
$italy.Range("A1:D11").Select()

# Make the new Spain sheet the active tab, with its own selection.
$newSheet.Activate()
$newSheet.Range("G8:G9").Select()
